# Update scripts with new TPM values.
# This mirrors a recomputation of the NATMI LR-pair stats (Gdf2-Acvr2a) after
# the TPM normalization changed. The "Receptor average/total expression value"
# (M2, N2) was rescaled, and all "derived specificity" columns (O, P, S, T)
# plus the edge weight columns (Q2, R2) were recomputed to reflect it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 19.539813446167
$ws.Range("R2").Value = 175.858321015503
$ws.Range("S2").Value = 0.3220556913988901
$ws.Range("T2").Value = 0.32205569139889

# Row 3
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("S3").Value = 0.5509544596378365
$ws.Range("T3").Value = 0.5509544596378364

# Row 4
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("S4").Value = 0.1269898489632735
$ws.Range("T4").Value = 0.1269898489632735
